$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 1.03
$ws.Range("N2").Value = 11
$ws.Range("O2").Value = 1.22
$ws.Range("P2").Value = 3.75
$ws.Range("Q2").Value = 1.93
$ws.Range("R2").Value = 1.93

# Row 5 updates
$ws.Range("Q5").Value = 1.88
$ws.Range("R5").Value = 1.98
